$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G2").Value = 1.8
$ws.Range("I2").Value = 4.4
$ws.Range("R2").Value = 5.4
$ws.Range("S2").Value = 7.3
$ws.Range("U2").Value = 14.5
$ws.Range("V2").Value = 17
$ws.Range("X2").Value = 7.3
$ws.Range("Z2").Value = 19
$ws.Range("AA2").Value = 120
$ws.Range("AB2").Value = 10
$ws.Range("AC2").Value = 23
$ws.Range("AD2").Value = 15
$ws.Range("AE2").Value = 75
$ws.Range("AF2").Value = 50
$ws.Range("AG2").Value = 65
$ws.Range("J3").Value = 1.4
$ws.Range("K3").Value = 3
$ws.Range("L3").Value = 2.2
$ws.Range("M3").Value = 1.67
$ws.Range("AI3").Value = 1.08
$ws.Range("N5").Value = 1.57
$ws.Range("O5").Value = 2.25
$ws.Range("P5").Value = 2.25
$ws.Range("Q5").Value = 1.57
$ws.Range("R5").Value = 8.5
$ws.Range("X5").Value = 7
$ws.Range("AA5").Value = 81
$ws.Range("G7").Value = 2.3
$ws.Range("I7").Value = 3.1
$ws.Range("Y7").Value = 6.5
$ws.Range("AJ7").Value = 8
$ws.Range("G8").Value = 1.65
$ws.Range("J8").Value = 1.5
$ws.Range("K8").Value = 2.63
$ws.Range("N8").Value = 1.53
$ws.Range("O8").Value = 2.38
$ws.Range("R8").Value = 5
$ws.Range("Y8").Value = 7.5
$ws.Range("Z8").Value = 23
$ws.Range("AD8").Value = 19
$ws.Range("AE8").Value = 67
$ws.Range("G10").Value = 2.35
$ws.Range("H10").Value = 3.3
$ws.Range("I10").Value = 3
$ws.Range("V10").Value = 21
$ws.Range("X10").Value = 9
$ws.Range("Z10").Value = 15
$ws.Range("G12").Value = 1.73
$ws.Range("I12").Value = 4.33
$ws.Range("U12").Value = 13
$ws.Range("V12").Value = 13
$ws.Range("Y12").Value = 7.5
$ws.Range("AB12").Value = 13
$ws.Range("AC12").Value = 23
$ws.Range("AD12").Value = 15
$ws.Range("AE12").Value = 51
$ws.Range("AH12").Value = 201
$ws.Range("I13").Value = 2.63
$ws.Range("P13").Value = 2.1
$ws.Range("Q13").Value = 1.67
$ws.Range("L15").Value = 2
$ws.Range("M15").Value = 1.85
$ws.Range("N16").Value = 1.33
$ws.Range("AI16").Value = 1.03
$ws.Range("AJ16").Value = 15
$ws.Range("G17").Value = 3.6
$ws.Range("N17").Value = 1.44
$ws.Range("O17").Value = 2.63
$ws.Range("S17").Value = 17
$ws.Range("AE17").Value = 19
$ws.Range("L18").Value = 2.1
$ws.Range("M18").Value = 1.7
$ws.Range("N18").Value = 1.44
$ws.Range("O18").Value = 2.63
$ws.Range("AI18").Value = 1.06
$ws.Range("AJ18").Value = 10
$ws.Range("G20").Value = 1.07
$ws.Range("H20").Value = 7.7
$ws.Range("I20").Value = 23
$ws.Range("L20").Value = 1.27
$ws.Range("M20").Value = 3.45
$ws.Range("P20").Value = 2.5
$ws.Range("Q20").Value = 1.47
$ws.Range("R20").Value = 9.5
$ws.Range("T20").Value = 11.25
$ws.Range("U20").Value = 5.3
$ws.Range("V20").Value = 10
$ws.Range("W20").Value = 35
$ws.Range("Y20").Value = 18.5
$ws.Range("Z20").Value = 40
$ws.Range("AA20").Value = 175
$ws.Range("AB20").Value = 70
$ws.Range("AC20").Value = 300
$ws.Range("AD20").Value = 80
$ws.Range("AF20").Value = 500
$ws.Range("AG20").Value = 250
$ws.Range("G21").Value = 5.2
$ws.Range("H21").Value = 3.85
$ws.Range("I21").Value = 1.52
$ws.Range("L21").Value = 1.53
$ws.Range("M21").Value = 2.18
$ws.Range("R21").Value = 15
$ws.Range("T21").Value = 13.5
$ws.Range("U21").Value = 75
$ws.Range("X21").Value = 14
$ws.Range("Y21").Value = 7
$ws.Range("Z21").Value = 11.75
$ws.Range("AB21").Value = 7.4
$ws.Range("AD21").Value = 6.8
$ws.Range("AE21").Value = 10
$ws.Range("AF21").Value = 9.5
$ws.Range("G22").Value = 1.67
$ws.Range("I22").Value = 4.5
$ws.Range("AI22").Value = 1.04
$ws.Range("AJ22").Value = 13
$ws.Range("G23").Value = 3.4
$ws.Range("I23").Value = 2
$ws.Range("L23").Value = 1.88
$ws.Range("M23").Value = 1.98
$ws.Range("N23").Value = 1.36
$ws.Range("O23").Value = 3
$ws.Range("T23").Value = 12
$ws.Range("Z23").Value = 13
$ws.Range("AC23").Value = 10
$ws.Range("AE23").Value = 19
$ws.Range("AF23").Value = 17
$ws.Range("L25").Value = 1.9
$ws.Range("M25").Value = 1.95
$ws.Range("J29").Value = 1.25
$ws.Range("K29").Value = 3.7
$ws.Range("L29").Value = 1.75
$ws.Range("M29").Value = 1.87
$ws.Range("P29").Value = 2.3
$ws.Range("Q29").Value = 1.57
$ws.Range("R29").Value = 5.1
$ws.Range("S29").Value = 4.9
$ws.Range("T29").Value = 7
$ws.Range("U29").Value = 6.5
$ws.Range("V29").Value = 9.25
$ws.Range("W29").Value = 24
$ws.Range("X29").Value = 9.25
$ws.Range("Z29").Value = 18.5
$ws.Range("AA29").Value = 90
$ws.Range("AC29").Value = 80
$ws.Range("AD29").Value = 30
$ws.Range("AE29").Value = 400
$ws.Range("AF29").Value = 150
$ws.Range("AG29").Value = 100
$ws.Range("I32").Value = 8
$ws.Range("S32").Value = 6.5
$ws.Range("T32").Value = 8.5
$ws.Range("Y32").Value = 8
$ws.Range("AA32").Value = 67
$ws.Range("AB32").Value = 17
$ws.Range("AC32").Value = 41
$ws.Range("AD32").Value = 23
$ws.Range("AI32").Value = 1.06
$ws.Range("AJ32").Value = 9.5
